$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "'45.638.49"
$ws.Range("E2").Value = "  -2.09%  "
$ws.Range("D3").Value = "'2.421.68"
$ws.Range("E3").Value = "  +5.51%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'300.44"
$ws.Range("E5").Value = "  -1.46%  "
$ws.Range("D6").Value = "'97.51"
$ws.Range("E6").Value = "  -3.47%  "
$ws.Range("D7").Value = "'0.565"
$ws.Range("E7").Value = "  -0.32%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "'0.512"
$ws.Range("E9").Value = "  -1.48%  "
$ws.Range("D10").Value = "'34.54"
$ws.Range("E10").Value = "  -5.96%  "
$ws.Range("D11").Value = "'0.0793"
$ws.Range("E11").Value = "  +0.27%  "
$ws.Range("D12").Value = "'7.19"
$ws.Range("E12").Value = "  -2.48%  "
$ws.Range("E13").Value = "  +0.55%  "
$ws.Range("D14").Value = "'2.789.78"
$ws.Range("E14").Value = "  +5.56%  "
$ws.Range("D15").Value = "'2.406.47"
$ws.Range("E15").Value = "  +5.22%  "
$ws.Range("D16").Value = "'14.18"
$ws.Range("E16").Value = "  +2.66%  "
$ws.Range("D17").Value = "'0.835"
$ws.Range("E17").Value = "  +3.08%  "
$ws.Range("D18").Value = "'45.632.83"
$ws.Range("E18").Value = "  -2.03%  "
$ws.Range("D19").Value = "'13.05"
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("D20").Value = "'0.0₃0954"
$ws.Range("E20").Value = "  +1.69%  "
$ws.Range("D21").Value = "'6.18"
$ws.Range("E21").Value = "  +2.52%  "
$ws.Range("D22").Value = "'67.43"
$ws.Range("E22").Value = "  +1.59%  "
$ws.Range("D23").Value = "'244.03"
$ws.Range("E23").Value = "  -1.51%  "
$ws.Range("D24").Value = "'2.81"
$ws.Range("E24").Value = "  -3.47%  "
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("D26").Value = "'1.94"
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("D27").Value = "'39.05"
$ws.Range("E27").Value = "  -9.63%  "
$ws.Range("D28").Value = "'2.22"
$ws.Range("E28").Value = "  -1.87%  "
$ws.Range("D29").Value = "'9.82"
$ws.Range("E29").Value = "  -0.40%  "
$ws.Range("D30").Value = "'3.87"
$ws.Range("E30").Value = "  +19.33%  "
$ws.Range("D31").Value = "'21.42"
$ws.Range("E31").Value = "  +7.33%  "
$ws.Range("D32").Value = "'5.61"
$ws.Range("E32").Value = "  -0.61%  "
$ws.Range("E33").Value = "  -2.78%  "
$ws.Range("D34").Value = "'147.66"
$ws.Range("E34").Value = "  +0.29%  "
$ws.Range("D35").Value = "'0.0778"
$ws.Range("E35").Value = "  -2.22%  "
$ws.Range("D36").Value = "'1.98"
$ws.Range("E36").Value = "  +11.43%  "
$ws.Range("D37").Value = "'0.113"
$ws.Range("E37").Value = "  -1.17%  "
$ws.Range("D38").Value = "'0.117"
$ws.Range("E38").Value = "  -0.73%  "
$ws.Range("D39").Value = "'15.41"
$ws.Range("E39").Value = "  -3.94%  "
$ws.Range("D40").Value = "'3.89"
$ws.Range("E40").Value = "  -4.44%  "
$ws.Range("E41").Value = "  -0.93%  "
$ws.Range("D42").Value = "'3.29"
$ws.Range("E42").Value = "  -2.59%  "
$ws.Range("D43").Value = "'1.964.14"
$ws.Range("E43").Value = "  +7.35%  "
$ws.Range("D45").Value = "'91.51"
$ws.Range("E45").Value = "  +4.65%  "
$ws.Range("E46").Value = "  -9.49%  "
$ws.Range("D47").Value = "'8.67"
$ws.Range("E47").Value = "  +10.47%  "
$ws.Range("D48").Value = "'100.30"
$ws.Range("E48").Value = "  +4.75%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "'2.658.25"
$ws.Range("E49").Value = "  +5.54%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "'0.186"
$ws.Range("E50").Value = "  -4.65%  "
$ws.Range("D51").Value = "'14.71"
$ws.Range("E51").Value = "  +9.15%  "
